$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing the old row 6 (attempt 4) down to row 7
$ws.Rows("6:6").Insert()

# Update row 5 (attempt 3) with new results from a re-run
$ws.Range("C5").Value = 89
$ws.Range("D5").Value = 2537
$ws.Range("G5").Value = "n/a"
$ws.Range("J5").Value = 0.7457
$ws.Range("J5").Font.Color = 1907741

# Fill in the newly inserted row 6 with duplicated data (copy of row5's prior loss/metric values,
# keeping the original row6 attempt info in columns A-F)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "neural network"
$ws.Range("C6").Value = 87
$ws.Range("D6").Value = 6289
$ws.Range("E6").Value = 8417
$ws.Range("F6").Value = "MSE"
$ws.Range("G6").Value = 9.0807
$ws.Range("H6").Value = "n/a"
$ws.Range("I6").Value = "n/a"
$ws.Range("J6").Value = 0.8324

# The attempt that got pushed down to row 7 is renumbered as attempt 5
$ws.Range("A7").Value = 5

# Update the selection and page setup, as left by the author after editing
[void]$ws.Range("K10").Select()
$ws.PageSetup.Orientation = 1
